$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy formatting (number format / font / style) from column E into the new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Remove the stray empty cells created in rows that never had data in columns C:K
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# Populate the new column D with the latest reporting period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1961800
$ws.Range("D9").Value = 704600
$ws.Range("D10").Value = 1257100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -78800
$ws.Range("D15").Value = 297900
$ws.Range("D17").Value = 1519400
$ws.Range("D18").Value = 442300
$ws.Range("D20").Value = -38900
$ws.Range("D21").Value = 701300
$ws.Range("D22").Value = 186200
$ws.Range("D23").Value = 217200
$ws.Range("D24").Value = 11700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 205500
$ws.Range("D27").Value = -275500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 300
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 38900
$ws.Range("D33").Value = -275200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -275200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 13600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 150400
$ws.Range("D44").Value = 22700
$ws.Range("D45").Value = 15400
$ws.Range("D46").Value = 202200
$ws.Range("D47").Value = 500
$ws.Range("D48").Value = 4288600
$ws.Range("D49").Value = 1770000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 87800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 6349100
$ws.Range("D57").Value = 143100
$ws.Range("D58").Value = 18500
$ws.Range("D59").Value = 161800
$ws.Range("D60").Value = 323400
$ws.Range("D61").Value = 3112000
$ws.Range("D62").Value = 92000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3527400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 1320300
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1501400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -275200
$ws.Range("D83").Value = 297900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 544200
$ws.Range("D91").Value = -457500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -153800
$ws.Range("D96").Value = -391400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -399900
$ws.Range("D101").Value = -1200
$ws.Range("D102").Value = -10600
